# Auto-generated Excel COM-interop script
# Applies value corrections to currentAveragePrice / Leve profit calculation
# columns (H-N) across all 8 sheets, as captured by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3024.875
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 3024.875
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 9074.625
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -9298.625

$ws.Range("H33").Value = 76.34999999999999
$ws.Range("I33").Value = 68.71429000000001
$ws.Range("K33").Value = 68.71429000000001
$ws.Range("M33").Value = 160.28571

$ws.Range("H94").Value = 2379.2
$ws.Range("I94").Value = 2474
$ws.Range("K94").Value = 2474
$ws.Range("M94").Value = -2023

$ws.Range("H138").Value = 4685.451
$ws.Range("I138").Value = 5538.8667
$ws.Range("J138").Value = 4494.388
$ws.Range("K138").Value = 16616.6001
$ws.Range("L138").Value = 13483.164
$ws.Range("M138").Value = -11476.6001
$ws.Range("N138").Value = -23763.164

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9915.741
$ws.Range("I32").Value = 8608.74
$ws.Range("J32").Value = 26253.25
$ws.Range("K32").Value = 8608.74
$ws.Range("L32").Value = 26253.25
$ws.Range("M32").Value = -8321.74
$ws.Range("N32").Value = -26827.25

$ws.Range("H110").Value = 866.75
$ws.Range("J110").Value = 1060
$ws.Range("L110").Value = 1060
$ws.Range("N110").Value = -5150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 5477.4
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 5477.4
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H94").Value = 770.96295
$ws.Range("I94").Value = 445.05
$ws.Range("K94").Value = 445.05
$ws.Range("M94").Value = 5.949999999999989

$ws.Range("H105").Value = 2250.8262
$ws.Range("I105").Value = 2162.5642
$ws.Range("J105").Value = 2742.5715
$ws.Range("K105").Value = 2162.5642
$ws.Range("L105").Value = 2742.5715
$ws.Range("M105").Value = -415.5641999999998
$ws.Range("N105").Value = -6236.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4832.8335
$ws.Range("I31").Value = 3666.6667
$ws.Range("J31").Value = 5999
$ws.Range("K31").Value = 3666.6667
$ws.Range("L31").Value = 5999
$ws.Range("M31").Value = -3371.6667
$ws.Range("N31").Value = -6589

$ws.Range("H34").Value = 4832.8335
$ws.Range("I34").Value = 3666.6667
$ws.Range("J34").Value = 5999
$ws.Range("K34").Value = 3666.6667
$ws.Range("L34").Value = 5999
$ws.Range("M34").Value = -3464.6667
$ws.Range("N34").Value = -6403

$ws.Range("H58").Value = 1978203.1
$ws.Range("I58").Value = 2719014.2
$ws.Range("K58").Value = 2719014.2
$ws.Range("M58").Value = -2718811.2

$ws.Range("H59").Value = 39284.668

$ws.Range("H62").Value = 4072.2856
$ws.Range("I62").Value = 4225
$ws.Range("J62").Value = 3868.6667
$ws.Range("K62").Value = 4225
$ws.Range("L62").Value = 3868.6667
$ws.Range("M62").Value = -3601
$ws.Range("N62").Value = -5116.6667

$ws.Range("H65").Value = 4072.2856
$ws.Range("I65").Value = 4225
$ws.Range("J65").Value = 3868.6667
$ws.Range("K65").Value = 21125
$ws.Range("L65").Value = 19343.3335
$ws.Range("M65").Value = -18005
$ws.Range("N65").Value = -25583.3335

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H99").Value = 11006
$ws.Range("I99").Value = 11006
$ws.Range("K99").Value = 11006
$ws.Range("M99").Value = -9508

$ws.Range("H126").Value = 11006
$ws.Range("I126").Value = 11006
$ws.Range("K126").Value = 33018
$ws.Range("M126").Value = -30548

$ws.Range("H131").Value = 58888
$ws.Range("J131").Value = 58888
$ws.Range("L131").Value = 58888
$ws.Range("N131").Value = -68968

$ws.Range("H132").Value = 2654.4666
$ws.Range("I132").Value = 1619.909
$ws.Range("K132").Value = 4859.727000000001
$ws.Range("M132").Value = -2329.727000000001

$ws.Range("H133").Value = 63899
$ws.Range("J133").Value = 63899
$ws.Range("L133").Value = 63899
$ws.Range("N133").Value = -68959

$ws.Range("H134").Value = 1329.0303
$ws.Range("I134").Value = 1169.76
$ws.Range("J134").Value = 1826.75
$ws.Range("K134").Value = 3509.28
$ws.Range("L134").Value = 5480.25
$ws.Range("M134").Value = -974.2799999999997
$ws.Range("N134").Value = -10550.25

$ws.Range("H136").Value = 1978203.1
$ws.Range("I136").Value = 2719014.2
$ws.Range("K136").Value = 8157042.600000001
$ws.Range("M136").Value = -8154492.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 452112.3
$ws.Range("I4").Value = 420365.8
$ws.Range("K4").Value = 1261097.4
$ws.Range("M4").Value = -1260985.4

$ws.Range("H112").Value = 42599.6
$ws.Range("J112").Value = 42599.6
$ws.Range("L112").Value = 127798.8
$ws.Range("N112").Value = -130014.8

$ws.Range("H113").Value = 75066.39999999999
$ws.Range("J113").Value = 1999.3846
$ws.Range("L113").Value = 5998.1538
$ws.Range("N113").Value = -10338.1538

$ws.Range("H131").Value = 11694.246
$ws.Range("J131").Value = 11694.246
$ws.Range("L131").Value = 35082.738
$ws.Range("N131").Value = -45162.738

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1795.7894
$ws.Range("I97").Value = 1997.1428
$ws.Range("K97").Value = 1997.1428
$ws.Range("M97").Value = -1501.1428

$ws.Range("H126").Value = 1827182.5
$ws.Range("J126").Value = 169319
$ws.Range("L126").Value = 507957
$ws.Range("N126").Value = -512897

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2199.4
$ws.Range("I68").Value = 1999.6666
$ws.Range("J68").Value = 2499
$ws.Range("K68").Value = 1999.6666
$ws.Range("L68").Value = 2499
$ws.Range("M68").Value = -1250.6666
$ws.Range("N68").Value = -3997

$ws.Range("H71").Value = 2199.4
$ws.Range("I71").Value = 1999.6666
$ws.Range("J71").Value = 2499
$ws.Range("K71").Value = 9998.333000000001
$ws.Range("L71").Value = 12495
$ws.Range("M71").Value = -6254.333000000001
$ws.Range("N71").Value = -19983

$ws.Range("H136").Value = 4880.25
$ws.Range("I136").Value = 3681
$ws.Range("J136").Value = 5599.8
$ws.Range("K136").Value = 11043
$ws.Range("L136").Value = 16799.4
$ws.Range("M136").Value = -8493
$ws.Range("N136").Value = -21899.4

$ws.Range("H139").Value = 69599
$ws.Range("J139").Value = 69599
$ws.Range("L139").Value = 69599
$ws.Range("N139").Value = -79879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1699.3334
$ws.Range("I81").Value = 1699.3334
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3398.6668
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2337.6668
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 1699.3334
$ws.Range("I84").Value = 1699.3334
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 16993.334
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -11689.334
$ws.Range("N84").ClearContents()

$ws.Range("H132").Value = 2039.8667
$ws.Range("I132").Value = 1059.6
$ws.Range("K132").Value = 3178.8
$ws.Range("M132").Value = -648.7999999999997

$ws.Range("H136").Value = 29243512
$ws.Range("I136").Value = 42737824
$ws.Range("K136").Value = 128213472
$ws.Range("M136").Value = -128210922

$ws.Range("H139").Value = 73919.60000000001
$ws.Range("J139").Value = 73919.60000000001
$ws.Range("L139").Value = 73919.60000000001
$ws.Range("N139").Value = -84199.60000000001
